$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Add logging mechanism..." task row. This shifts every
# subsequent task up by one row (so the old row 12 disappears and the
# used range becomes A1:B11).
$ws.Rows(2).Delete()

# The row that used to hold "Error handling..." (value 21) is now row 2.
# Overwrite it with the new diagnostics task that was added as part of
# this "first raft of diagnostics" commit.
$ws.Range("A2").Value = "Make exporter a GUP.  Build UI & hook data into the max files"
$ws.Range("B2").Value = 5

# Match the author's final cursor position.
$ws.Range("B3").Select() | Out-Null
